# edit.ps1
# Implements the commit "feat: add 2022-Q3 data":
#   1. Insert a new worksheet "2022-Q3" right after "总计" (duplicated from
#      the existing "2022-Q2" sheet so it inherits identical formatting),
#      which pushes every later quarterly sheet down by one tab position
#      (their content is otherwise unchanged).
#   2. Populate "2022-Q3" with the fund-holding table for that quarter.
#   3. Update the "总计" (totals) sheet: insert a new leading data row for
#      2022-Q3 (count=32, value=11.84) and push the previously existing
#      rows down by one, re-numbering the running index column (A).
#
# NOTE: figures such as fund codes / percentages are stored as literal
# TEXT in this workbook (matching every other quarterly sheet), not as
# numbers - Excel would otherwise strip leading zeros / reformat them, so
# those are written with a leading apostrophe to force text entry, exactly
# as a user typing into the grid would do.

$wb = $excel.ActiveWorkbook

function Set-TextCell {
    param($cell, [string]$val)
    $cell.Value = "'" + $val
}

# ---------------------------------------------------------------------------
# Step 1: duplicate "2022-Q2" (so the new tab inherits the same column
# widths / header style / borders) and move the copy right after "总计".
# ---------------------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item("2022-Q2")
$totalSheet = $wb.Worksheets.Item("总计")
$templateSheet.Copy([System.Reflection.Missing]::Value, $totalSheet)

$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# Step 2: fill in the 2022-Q3 fund table.
# Columns: A=index, B=基金代码, C=基金名称, D=基金规模, E=股票总仓位,
#          F=仓位占比, G=持有市值(亿元), H=仓位排名
# ---------------------------------------------------------------------------
$q3Data = @(
    ,@(0, "000979", "景顺长城沪港深精选股票", "20.32", "80.04", "7.99", "1.6236", 1, 4)
    ,@(1, "008850", "景顺长城价值稳进三年定期开放灵活配置混合", "17.29", "84.84", "8.84", "1.5284", 1, 2)
    ,@(2, "011081", "国投瑞银港股通混合C", "19.35", "90.31", "7.65", "1.4803", 1, 4)
    ,@(3, "007110", "国投瑞银港股通价值发现混合", "19.35", "90.31", "7.65", "1.4803", 1, 4)
    ,@(4, "100061", "富国中国中小盘混合（QDII）人民币", "35.11", "83.32", "2.66", "0.9339", 1, 8)
    ,@(5, "010591", "富国中国中小盘混合（QDII）美元", "35.11", "83.32", "2.66", "0.9339", 1, 8)
    ,@(6, "008715", "景顺长城价值驱动一年持有期灵活配置混合", "8.44", "90.91", "9.54", "0.8052", 1, 5)
    ,@(7, "009098", "景顺长城价值领航两年持有期混合", "7.16", "90.92", "9.50", "0.6802", 1, 3)
    ,@(8, "008060", "景顺长城价值边际灵活配置混合A", "5.45", "80.34", "7.89", "0.4300", 1, 6)
    ,@(9, "009846", "富兰克林国海港股通远见价值混合", "12.78", "83.81", "3.10", "0.3962", 1, 5)
    ,@(10, "010010", "国投瑞银港股通6个月定期开放股票", "6.96", "90.21", "5.49", "0.3821", 1, 6)
    ,@(11, "009983", "永赢港股通品质生活慧选混合", "9.17", "60.68", "3.07", "0.2815", 1, 10)
    ,@(12, "012640", "鹏华稳健鸿利一年持有期混合A", "2.61", "92.98", "9.81", "0.2560", 1, 2)
    ,@(13, "008134", "鹏华优选价值股票", "1.80", "92.72", "9.41", "0.1694", 1, 2)
    ,@(14, "011574", "鹏华领航一年持有期混合A", "1.20", "92.84", "9.48", "0.1138", 1, 1)
    ,@(15, "011575", "鹏华领航一年持有期混合C", "0.91", "92.84", "9.48", "0.0863", 1, 1)
    ,@(16, "015779", "景顺长城价值边际灵活配置混合C", "0.79", "80.34", "7.89", "0.0623", 1, 6)
    ,@(17, "007107", "太平 MSCI 香港价值增强指数A", "0.97", "90.37", "3.96", "0.0384", 1, 7)
    ,@(18, "005052", "上投摩根标普港股通低波红利指数C", "1.37", "92.94", "2.30", "0.0315", 1, 8)
    ,@(19, "005051", "上投摩根标普港股通低波红利指数A", "1.36", "92.94", "2.30", "0.0313", 1, 8)
    ,@(20, "004266", "招商沪港深科技创新主题精选灵活配置混合A", "0.92", "90.52", "2.92", "0.0269", 1, 6)
    ,@(21, "010783", "德邦沪港深龙头混合A", "0.55", "84.96", "4.51", "0.0248", 1, 5)
    ,@(22, "010784", "德邦沪港深龙头混合C", "0.36", "84.96", "4.51", "0.0162", 1, 5)
    ,@(23, "012641", "鹏华稳健鸿利一年持有期混合C", "0.10", "92.98", "9.81", "0.0098", 1, 2)
    ,@(24, "010754", "招商沪港深科技创新主题精选灵活配置混合C", "0.25", "90.52", "2.92", "0.0073", 1, 6)
    ,@(25, "011647", "博时港股通红利精选混合A", "0.11", "82.44", "3.13", "0.0034", 1, 9)
    ,@(26, "501303", "广发恒生中型股指数（LOF）A", "0.21", "89.12", "1.54", "0.0032", 1, 8)
    ,@(27, "004996", "广发恒生中型股指数（LOF）C", "0.09", "89.12", "1.54", "0.0014", 1, 8)
    ,@(28, "160922", "大成恒生综合中小型股指数（QDII-LOF）A", "0.09", "86.62", "1.12", "0.0010", 1, 8)
    ,@(29, "011648", "博时港股通红利精选混合C", "0.02", "82.44", "3.13", "0.0006", 1, 9)
    ,@(30, "008972", "大成恒生综合中小型股指数C", "0.02", "86.62", "1.12", "0.0002", 1, 8)
    ,@(31, "007108", "太平 MSCI 香港价值增强指数C", "0.00", "90.37", "3.96", 0, 0, 7)
)


$lastRow = 1 + $q3Data.Count   # header row + data rows

# Extend column A's header-row style down through every data row (the
# template sheet only had 17 data rows; this quarter has 32).
$q3Sheet.Range("A18").Copy()
if ($lastRow -gt 18) {
    $q3Sheet.Range("A19:A" + $lastRow).PasteSpecial(-4122)
}

foreach ($row in $q3Data) {
    $r = [int]$row[0] + 2
    $q3Sheet.Cells.Item($r, 1).Value = [int]$row[0]
    Set-TextCell $q3Sheet.Cells.Item($r, 2) $row[1]
    $q3Sheet.Cells.Item($r, 3).Value = $row[2]
    Set-TextCell $q3Sheet.Cells.Item($r, 4) $row[3]
    Set-TextCell $q3Sheet.Cells.Item($r, 5) $row[4]
    Set-TextCell $q3Sheet.Cells.Item($r, 6) $row[5]
    if ([int]$row[7] -eq 1) {
        Set-TextCell $q3Sheet.Cells.Item($r, 7) $row[6]
    } else {
        $q3Sheet.Cells.Item($r, 7).Value = $row[6]
    }
    $q3Sheet.Cells.Item($r, 8).Value = [int]$row[8]
}

# ---------------------------------------------------------------------------
# Step 3: update the "总计" summary sheet - insert the 2022-Q3 row at the
# top of the data and push the rest down by one row.
# Columns: A=index, B=日期, C=持有数量(只), D=持有市值(亿元)
# ---------------------------------------------------------------------------
$totalsData = @(
    ,@(0, "2022-Q3", 32, 11.84)
    ,@(1, "2022-Q2", 17, 8.16)
    ,@(2, "2022-Q1", 26, 12.05)
    ,@(3, "2021-Q4", 20, 8.96)
    ,@(4, "2021-Q3", 15, 4.62)
    ,@(5, "2021-Q2", 14, 5.21)
    ,@(6, "2021-Q1", 22, 6.91)
    ,@(7, "2020-Q4", 6, 1.08)
)


# Row 9 is brand new - give column A the same style as the previous last
# row (A8) before filling in values.
$totalSheet.Range("A8").Copy()
$totalSheet.Range("A9").PasteSpecial(-4122)

foreach ($row in $totalsData) {
    $r = [int]$row[0] + 2
    $totalSheet.Cells.Item($r, 1).Value = [int]$row[0]
    $totalSheet.Cells.Item($r, 2).Value = $row[1]
    $totalSheet.Cells.Item($r, 3).Value = [int]$row[2]
    $totalSheet.Cells.Item($r, 4).Value = $row[3]
}
